# Add MAE column between MSE/R2 and Tipo, and refresh updated prediction
# values (commit: "Creado grafico de tipos de modelo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D ("Tipo") to make room
# for the new "MAE" column; this shifts "Tipo" to column E.
$ws.Range("D1").EntireColumn.Insert()

# Copy the header style from the old header cells (now at A1:C1) onto the
# new D1 header cell, then set its value / style to match C1 (header look).
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "MAE"

# Header row values
$ws.Range("A1").Value = "Enfermedad"
$ws.Range("B1").Value = "MSE"
$ws.Range("C1").Value = "R2"
$ws.Range("D1").Value = "MAE"
$ws.Range("E1").Value = "Tipo"

# Updated MSE (B), R2 (C) values + new MAE (D) values, per row.
$data = @(
    @{ Row = 2;  B = 0.4101026717820584; C = 0.9331900416269556; D = 0.5013618022923331 },
    @{ Row = 3;  B = 2.878893579773742;  C = 0.9588727023492972; D = 1.266131410638987 },
    @{ Row = 4;  B = 1.093677535067957;  C = 0.9459767614280777; D = 0.8047074834725575 },
    @{ Row = 5;  B = 1.56219342376502;   C = 0.9989738867125377; D = 0.9389581486305347 },
    @{ Row = 6;  B = 1.191061667365153;  C = 0.9870213363411285; D = 0.9111527565454832 },
    @{ Row = 7;  B = 1.271977756333023;  C = 0.9992868418203459; D = 0.834472912226307 },
    @{ Row = 8;  B = 1.439159786975202;  C = 0.9983515097052726; D = 0.9350345842168691 },
    @{ Row = 9;  B = 4.877178855194897;  C = 0.94169669696305;   D = 1.677266164849096 },
    @{ Row = 10; B = 0.8466835870385386; C = 0.9974760624516991; D = 0.7276269921692753 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
}
